# Applies the Chocobo_Profits scheduled-runner update: refreshes the
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N)
# for the rows whose underlying market data changed, across all 8 job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).ClearContents()
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).ClearContents()
$ws.Cells.Item(70, 8).Value = 3451.7334
$ws.Cells.Item(70, 9).Value = 2002
$ws.Cells.Item(70, 10).Value = 3555.2856
$ws.Cells.Item(70, 11).Value = 6006
$ws.Cells.Item(70, 12).Value = 10665.8568
$ws.Cells.Item(70, 13).Value = -5736
$ws.Cells.Item(70, 14).Value = -11205.8568
$ws.Cells.Item(73, 8).Value = 3451.7334
$ws.Cells.Item(73, 9).Value = 2002
$ws.Cells.Item(73, 10).Value = 3555.2856
$ws.Cells.Item(73, 11).Value = 6006
$ws.Cells.Item(73, 12).Value = 10665.8568
$ws.Cells.Item(73, 13).Value = -5070
$ws.Cells.Item(73, 14).Value = -12537.8568
$ws.Cells.Item(112, 8).Value = 1318.0944
$ws.Cells.Item(112, 9).Value = 831.2857
$ws.Cells.Item(112, 10).Value = 1392.174
$ws.Cells.Item(112, 11).Value = 2493.8571
$ws.Cells.Item(112, 12).Value = 4176.522
$ws.Cells.Item(112, 13).Value = -1385.8571
$ws.Cells.Item(112, 14).Value = -6392.522
$ws.Cells.Item(118, 8).Value = 739.6923
$ws.Cells.Item(118, 9).Value = 606
$ws.Cells.Item(118, 10).Value = 953.6
$ws.Cells.Item(118, 11).Value = 1818
$ws.Cells.Item(118, 12).Value = 2860.8
$ws.Cells.Item(118, 13).Value = -161
$ws.Cells.Item(118, 14).Value = -6174.8
$ws.Cells.Item(132, 8).Value = 242022.08
$ws.Cells.Item(132, 9).Value = 4053.3333
$ws.Cells.Item(132, 10).Value = 1114574.1
$ws.Cells.Item(132, 11).Value = 12159.9999
$ws.Cells.Item(132, 12).Value = 3343722.3
$ws.Cells.Item(132, 13).Value = -9629.999899999999
$ws.Cells.Item(132, 14).Value = -3348782.3
$ws.Cells.Item(135, 8).Value = 221.85715
$ws.Cells.Item(135, 9).Value = 193.03703
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 11).Value = 1737.33327
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 13).Value = 797.6667300000001
$ws.Cells.Item(135, 14).Value = -14070
$ws.Cells.Item(137, 8).Value = 3408.7693
$ws.Cells.Item(137, 9).Value = 2162.2727
$ws.Cells.Item(137, 10).Value = 4322.8667
$ws.Cells.Item(137, 11).Value = 6486.8181
$ws.Cells.Item(137, 12).Value = 12968.6001
$ws.Cells.Item(137, 13).Value = -3936.8181
$ws.Cells.Item(138, 8).Value = 3726.87
$ws.Cells.Item(138, 9).Value = 865.13336
$ws.Cells.Item(138, 10).Value = 4953.3286
$ws.Cells.Item(138, 11).Value = 2595.40008
$ws.Cells.Item(138, 12).Value = 14859.9858
$ws.Cells.Item(138, 13).Value = 2544.59992
$ws.Cells.Item(138, 14).Value = -25139.9858

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5579.5693
$ws.Cells.Item(32, 9).Value = 4455.231
$ws.Cells.Item(32, 10).Value = 10076.923
$ws.Cells.Item(32, 11).Value = 4455.231
$ws.Cells.Item(32, 12).Value = 10076.923
$ws.Cells.Item(32, 13).Value = -4168.231
$ws.Cells.Item(32, 14).Value = -10650.923
$ws.Cells.Item(43, 8).Value = 17444
$ws.Cells.Item(43, 9).Value = 5000
$ws.Cells.Item(43, 10).Value = 29888
$ws.Cells.Item(43, 11).Value = 5000
$ws.Cells.Item(43, 12).Value = 29888
$ws.Cells.Item(43, 13).Value = -4687
$ws.Cells.Item(43, 14).Value = -30514
$ws.Cells.Item(61, 8).Value = 970
$ws.Cells.Item(61, 9).Value = 660.9524
$ws.Cells.Item(61, 10).Value = 3133.3333
$ws.Cells.Item(61, 11).Value = 660.9524
$ws.Cells.Item(61, 12).Value = 3133.3333
$ws.Cells.Item(61, 13).Value = -448.9524
$ws.Cells.Item(61, 14).Value = -3557.3333
$ws.Cells.Item(109, 8).Value = 30490.477
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 30490.477
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 30490.477
$ws.Cells.Item(109, 14).Value = -33264.477
$ws.Cells.Item(132, 8).Value = 2000.25
$ws.Cells.Item(132, 9).Value = 881.12
$ws.Cells.Item(132, 10).Value = 4543.727
$ws.Cells.Item(132, 11).Value = 2643.36
$ws.Cells.Item(132, 12).Value = 13631.181
$ws.Cells.Item(132, 13).Value = -113.3600000000001
$ws.Cells.Item(132, 14).Value = -18691.181
$ws.Cells.Item(136, 8).Value = 970
$ws.Cells.Item(136, 9).Value = 660.9524
$ws.Cells.Item(136, 10).Value = 3133.3333
$ws.Cells.Item(136, 11).Value = 1982.8572
$ws.Cells.Item(136, 12).Value = 9399.999899999999
$ws.Cells.Item(136, 13).Value = 567.1428000000001
$ws.Cells.Item(136, 14).Value = -14499.9999

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 9).Value = 1466.5
$ws.Cells.Item(86, 10).Value = 2250.2856
$ws.Cells.Item(86, 11).Value = 1466.5
$ws.Cells.Item(86, 12).Value = 2250.2856
$ws.Cells.Item(86, 13).Value = -343.5
$ws.Cells.Item(86, 14).Value = -4496.2856
$ws.Cells.Item(89, 9).Value = 1466.5
$ws.Cells.Item(89, 10).Value = 2250.2856
$ws.Cells.Item(89, 11).Value = 7332.5
$ws.Cells.Item(89, 12).Value = 11251.428
$ws.Cells.Item(89, 13).Value = -1716.5
$ws.Cells.Item(89, 14).Value = -22483.428
$ws.Cells.Item(134, 8).Value = 1480.72
$ws.Cells.Item(134, 9).Value = 997.09753
$ws.Cells.Item(134, 10).Value = 3683.889
$ws.Cells.Item(134, 11).Value = 2991.29259
$ws.Cells.Item(134, 12).Value = 11051.667
$ws.Cells.Item(134, 13).Value = -456.29259
$ws.Cells.Item(134, 14).Value = -16121.667

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2018.6716
$ws.Cells.Item(58, 9).Value = 1771.193
$ws.Cells.Item(58, 10).Value = 3429.3
$ws.Cells.Item(58, 11).Value = 1771.193
$ws.Cells.Item(58, 12).Value = 3429.3
$ws.Cells.Item(58, 13).Value = -1568.193
$ws.Cells.Item(58, 14).Value = -3835.3
$ws.Cells.Item(132, 8).Value = 2665.98
$ws.Cells.Item(132, 9).Value = 2021.4762
$ws.Cells.Item(132, 10).Value = 6049.625
$ws.Cells.Item(132, 11).Value = 6064.4286
$ws.Cells.Item(132, 12).Value = 18148.875
$ws.Cells.Item(132, 13).Value = -3534.4286
$ws.Cells.Item(132, 14).Value = -23208.875
$ws.Cells.Item(136, 8).Value = 2018.6716
$ws.Cells.Item(136, 9).Value = 1771.193
$ws.Cells.Item(136, 10).Value = 3429.3
$ws.Cells.Item(136, 11).Value = 5313.579
$ws.Cells.Item(136, 12).Value = 10287.9
$ws.Cells.Item(136, 13).Value = -2763.579
$ws.Cells.Item(136, 14).Value = -15387.9

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 12358.909
$ws.Cells.Item(26, 9).Value = 40183.332
$ws.Cells.Item(26, 10).Value = 1924.75
$ws.Cells.Item(26, 11).Value = 120549.996
$ws.Cells.Item(26, 12).Value = 5774.25
$ws.Cells.Item(26, 13).Value = -120261.996
$ws.Cells.Item(26, 14).Value = -6350.25
$ws.Cells.Item(113, 8).Value = 787.1070999999999
$ws.Cells.Item(113, 9).Value = 619.4706
$ws.Cells.Item(113, 10).Value = 1046.1818
$ws.Cells.Item(113, 11).Value = 1858.4118
$ws.Cells.Item(113, 12).Value = 3138.5454
$ws.Cells.Item(113, 13).Value = 311.5882000000001
$ws.Cells.Item(113, 14).Value = -7478.5454
$ws.Cells.Item(136, 8).Value = 3231.3914
$ws.Cells.Item(136, 9).Value = 3130
$ws.Cells.Item(136, 10).Value = 3363.2
$ws.Cells.Item(136, 11).Value = 9390
$ws.Cells.Item(136, 12).Value = 10089.6
$ws.Cells.Item(136, 13).Value = -4290
$ws.Cells.Item(136, 14).Value = -20289.6
$ws.Cells.Item(137, 8).Value = 3047.0454
$ws.Cells.Item(137, 9).Value = 2000
$ws.Cells.Item(137, 10).Value = 3096.9048
$ws.Cells.Item(137, 11).Value = 6000
$ws.Cells.Item(137, 12).Value = 9290.714399999999
$ws.Cells.Item(137, 13).Value = -900
$ws.Cells.Item(137, 14).Value = -19490.7144

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4155.2856
$ws.Cells.Item(132, 9).Value = 2339.5715
$ws.Cells.Item(132, 10).Value = 5971
$ws.Cells.Item(132, 11).Value = 7018.7145
$ws.Cells.Item(132, 12).Value = 17913
$ws.Cells.Item(132, 13).Value = -4488.7145

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1275.9
$ws.Cells.Item(16, 9).Value = 1325.2858
$ws.Cells.Item(16, 10).Value = 1160.6666
$ws.Cells.Item(16, 11).Value = 1325.2858
$ws.Cells.Item(16, 12).Value = 1160.6666
$ws.Cells.Item(16, 13).Value = -1155.2858
$ws.Cells.Item(16, 14).Value = -1500.6666
$ws.Cells.Item(68, 8).Value = 777.2787
$ws.Cells.Item(68, 9).Value = 662.3103599999999
$ws.Cells.Item(68, 10).Value = 3000
$ws.Cells.Item(68, 11).Value = 662.3103599999999
$ws.Cells.Item(68, 12).Value = 3000
$ws.Cells.Item(68, 13).Value = 86.68964000000005
$ws.Cells.Item(68, 14).Value = -4498
$ws.Cells.Item(71, 8).Value = 777.2787
$ws.Cells.Item(71, 9).Value = 662.3103599999999
$ws.Cells.Item(71, 10).Value = 3000
$ws.Cells.Item(71, 11).Value = 3311.5518
$ws.Cells.Item(71, 12).Value = 15000
$ws.Cells.Item(71, 13).Value = 432.4482000000003
$ws.Cells.Item(71, 14).Value = -22488
$ws.Cells.Item(82, 8).Value = 1443.2222
$ws.Cells.Item(82, 9).Value = 642.5454999999999
$ws.Cells.Item(82, 10).Value = 1993.6875
$ws.Cells.Item(82, 11).Value = 642.5454999999999
$ws.Cells.Item(82, 12).Value = 1993.6875
$ws.Cells.Item(82, 13).Value = -281.5454999999999
$ws.Cells.Item(82, 14).Value = -2715.6875
$ws.Cells.Item(85, 8).Value = 1443.2222
$ws.Cells.Item(85, 9).Value = 642.5454999999999
$ws.Cells.Item(85, 10).Value = 1993.6875
$ws.Cells.Item(85, 11).Value = 642.5454999999999
$ws.Cells.Item(85, 12).Value = 1993.6875
$ws.Cells.Item(85, 13).Value = 605.4545000000001
$ws.Cells.Item(85, 14).Value = -4489.6875
$ws.Cells.Item(100, 8).Value = 1852.7368
$ws.Cells.Item(100, 9).Value = 1649.8572
$ws.Cells.Item(100, 10).Value = 2420.8
$ws.Cells.Item(100, 11).Value = 1649.8572
$ws.Cells.Item(100, 12).Value = 2420.8
$ws.Cells.Item(100, 13).Value = -1108.8572
$ws.Cells.Item(132, 8).Value = 3175.56
$ws.Cells.Item(132, 9).Value = 1360.4166
$ws.Cells.Item(132, 10).Value = 7843.0713
$ws.Cells.Item(132, 11).Value = 4081.2498
$ws.Cells.Item(132, 12).Value = 23529.2139
$ws.Cells.Item(132, 13).Value = -1551.2498
$ws.Cells.Item(132, 14).Value = -28589.2139

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 352.86957
$ws.Cells.Item(113, 9).Value = 267.4
$ws.Cells.Item(113, 10).Value = 418.6154
$ws.Cells.Item(113, 11).Value = 802.1999999999999
$ws.Cells.Item(113, 12).Value = 1255.8462
$ws.Cells.Item(113, 13).Value = 1367.8
$ws.Cells.Item(113, 14).Value = -5595.8462
$ws.Cells.Item(132, 8).Value = 7753722.5
$ws.Cells.Item(132, 9).Value = 1018.65515
$ws.Cells.Item(132, 10).Value = 23812896
$ws.Cells.Item(132, 11).Value = 3055.96545
$ws.Cells.Item(132, 12).Value = 71438688
$ws.Cells.Item(132, 13).Value = -525.9654500000001
$ws.Cells.Item(132, 14).Value = -71443748
$ws.Cells.Item(136, 8).Value = 2583.5881
$ws.Cells.Item(136, 9).Value = 721.2692
$ws.Cells.Item(136, 10).Value = 8636.125
$ws.Cells.Item(136, 11).Value = 2163.8076
$ws.Cells.Item(136, 12).Value = 25908.375
$ws.Cells.Item(136, 13).Value = 386.1923999999999
$ws.Cells.Item(136, 14).Value = -31008.375
